# "Script now functions without upper or lower limit to ISO ppm"
# Update the test boundary values on Sheet1 so the lower bound (A1) goes
# negative and the upper bound (A3) goes above its previous ceiling,
# then leave the selection where the user last edited (A3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = -1.02
$ws.Range("A3").Value = 15.05

$ws.Activate()
$ws.Range("A3").Select()
